$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (shared string) renames
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# Column C (GDP) value updates
$ws.Range("C2").Value = 2870.311589353206
$ws.Range("C3").Value = 5596.139681459835
$ws.Range("C4").Value = 1460.056109840828
$ws.Range("C5").Value = 4729.735976516416
$ws.Range("C6").Value = 3972.630273980753
$ws.Range("C7").Value = 6128.19547247793
$ws.Range("C8").Value = 10594.98659239237
$ws.Range("C9").Value = 1909.084588129339
$ws.Range("C10").Value = 9502.243585046588
$ws.Range("C11").Value = 2100.656463590606
$ws.Range("C12").Value = 3928.450391496945
$ws.Range("C13").Value = 743.403784726004
$ws.Range("C14").Value = 2812.435974421079
$ws.Range("C15").Value = 665.6274194933962
$ws.Range("C16").Value = 1503.870423231357
$ws.Range("C17").Value = 10385.96443195552
$ws.Range("C18").Value = 1955.461557360978
$ws.Range("C19").Value = 6336.709213679884
$ws.Range("C20").Value = 4355.934938677345
$ws.Range("C21").Value = 5082.354756663512
$ws.Range("C23").Value = 777.227218443918
$ws.Range("C24").Value = 5885.254624554112
$ws.Range("C25").Value = 2965.153206179127
$ws.Range("C26").Value = 1577.487171555845
$ws.Range("C27").Value = 10883.31535948899
$ws.Range("C28").Value = 2024.117324382548
$ws.Range("C29").Value = 11627.81065059172
$ws.Range("C30").Value = 6711.616186806423
$ws.Range("C31").Value = 4479.398934239905
$ws.Range("C32").Value = 5360.226632400601
$ws.Range("C33").Value = 4209.874800894355
$ws.Range("C34").Value = 2860.874335573629
$ws.Range("C35").Value = 5642.578115155247
$ws.Range("C36").Value = 1657.651524528445
$ws.Range("C37").Value = 6911.59200404802
$ws.Range("C38").Value = 2094.024217383061
$ws.Range("C39").Value = 4394.543881413723
$ws.Range("C40").Value = 2887.250212489506
$ws.Range("C41").Value = 5919.20956823756
$ws.Range("C42").Value = 1716.389195271215
$ws.Range("C43").Value = 7200.731056811853
$ws.Range("C44").Value = 2201.396847776877
$ws.Range("C45").Value = 4699.493713911862
$ws.Range("C46").Value = 3008.669179463094
$ws.Range("C47").Value = 2612.856880840196
$ws.Range("C48").Value = 1775.027517189621
$ws.Range("C49").Value = 5996.49696468919
$ws.Range("C50").Value = 3012.536723186288
$ws.Range("C51").Value = 2735.187532014817
$ws.Range("C52").Value = 1836.014008604312
$ws.Range("C53").Value = 6114.227214287786
$ws.Range("C54").Value = 10239.48134799327
$ws.Range("C55").Value = 7633.969039669125
$ws.Range("C56").Value = 2854.757682901436
$ws.Range("C57").Value = 1895.214690888655
$ws.Range("C58").Value = 6262.368904654469

# AL15 flag flip 0 -> 1
$ws.Range("AL15").Value = 1
